$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.461.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.85"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.40"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +15.54%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.290"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0668"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0994"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.060.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.804.19"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.628"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.429.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0766"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.64"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0511"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.82"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.315.37"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.93"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +14.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.86%  "
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("E41").Value = "  +6.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.938"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0519"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.960.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.68"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("E51").Value = "  +1.46%  "
